# Edit sheet Card24 by admin
# Remove the trailing blank row (row 13) from the Card24 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")
$ws.Rows.Item(13).Delete()
